# Applies the "Added room locations functionalities" edit to the ACPH
# test-report template: refreshes the client/AHU header block, the
# instrument-calibration block, and the per-room obtained-results table
# with new sample data.
#
# Each target value lives in its own dedicated paragraph (confirmed via
# $d.Paragraphs enumeration), so we address paragraphs by their stable
# document-level index and overwrite the paragraph Range's Text. That
# keeps every run's formatting (rFonts/color/sz/shd) untouched and only
# rewrites the <w:t> content, exactly like the source diff does run-by-run.

$d = $word.ActiveDocument

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Text = $newText
}

# --- Client / Ref no / Area of Test block ---
Set-ParaText 13 "7790"                        # Sun Microsystem Pvt Ltd -> 7790
Set-ParaText 14 "7790"                        # 77908954 -> 7790
Set-ParaText 16 "12345678"                    # Chennai -> 12345678

# --- Plant / Date of Test / AHU Number block ---
Set-ParaText 19 "7790"                        # Cuddalore -> 7790
Set-ParaText 21 "15-11-2023 12.00.00 AM"      # 18-11-2023 -> 15-11-2023

# --- Instrument Used / Instrument Serial Number / Calibrated on block ---
Set-ParaText 29 "Ammenometer"                 # adada -> Ammenometer
Set-ParaText 30 "89991"                       # adad -> 89991
Set-ParaText 31 "16-11-2023 12.00.00 AM"      # 18-11-2023 -> 16-11-2023

# --- Make Model / Calibration Due on block ---
# This paragraph holds two independent runs ("adad" + " " + "adad"); both
# runs individually become "12345678" while the space run is untouched.
$p34 = $d.Paragraphs.Item(34)
$firstWord = $p34.Range.Words.Item(1)
$firstRange = $d.Range($firstWord.Start, $firstWord.End - 1)
$firstRange.Text = "12345678"

$p34 = $d.Paragraphs.Item(34)
$secondWord = $p34.Range.Words.Item(2)
$secondRange = $d.Range($secondWord.Start, $secondWord.End)
$secondRange.Text = "12345678"

Set-ParaText 35 "24-11-2023 12.00.00 AM"      # 07-11-2023 -> 24-11-2023

# --- Obtained Test Results table: room row "RSVD-2" / S3 ---
Set-ParaText 76 "RSVD-4"                      # RSVD-2 -> RSVD-4
Set-ParaText 80 "L1"                          # S3 -> L1
Set-ParaText 81 "0.98"                        # 0.25 -> 0.98
Set-ParaText 82 "189"                         # 1896 -> 189
Set-ParaText 83 "17"                          # 1856 -> 17
Set-ParaText 84 "717"                         # 4856 -> 717
Set-ParaText 85 "7187"                        # 5698 -> 7187
Set-ParaText 86 "1718"                        # 9874 -> 1718
Set-ParaText 87 "1966"                        # 4836 -> 1966
Set-ParaText 88 "1927"                        # 1209 -> 1927
Set-ParaText 91 "3384"                        # 4670 -> 3384
Set-ParaText 94 "2078"                        # 2500 -> 2078
Set-ParaText 97 "98"                          # 112 -> 98

# --- Obtained Test Results table: room row S4 ---
Set-ParaText 101 "L2"                         # S4 -> L2
Set-ParaText 102 "0.98"                       # 0.56 -> 0.98
Set-ParaText 103 "7187"                       # 1890 -> 7187
Set-ParaText 104 "87"                         # 9009 -> 87
Set-ParaText 105 "78"                         # 1909 -> 78
Set-ParaText 106 "78"                         # 9010 -> 78
Set-ParaText 107 "7"                          # 9080 -> 7
Set-ParaText 108 "1487"                       # 6180 -> 1487
Set-ParaText 109 "1457"                       # 3461 -> 1457
